# Update the "Förändrad" (changed) date in column C for rows 2-8
# from 2023-10-05 (45204) to 2023-10-08 (45207).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 45207
}
